$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-10-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-22 Tuesday", 2) | Out-Null

# Update the division problems in the table, cell by cell (row, col) so
# that cells whose new value equals another cell's old value are not
# accidentally clobbered by a global find/replace pass.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "96÷7="

$cell = $t.Cell(1, 2)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "99÷7="

$cell = $t.Cell(1, 3)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "82÷2="

$cell = $t.Cell(1, 4)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "95÷4="

$cell = $t.Cell(1, 5)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "24÷9="

$cell = $t.Cell(5, 1)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "77÷9="

$cell = $t.Cell(5, 2)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "86÷9="

$cell = $t.Cell(5, 3)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "42÷2="

$cell = $t.Cell(5, 4)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "90÷7="

$cell = $t.Cell(5, 5)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "14÷8="

$cell = $t.Cell(9, 1)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "18÷6="

$cell = $t.Cell(9, 2)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "88÷2="

$cell = $t.Cell(9, 3)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "58÷2="

$cell = $t.Cell(9, 4)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "37÷8="

$cell = $t.Cell(9, 5)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "78÷6="

$cell = $t.Cell(13, 1)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "72÷9="

$cell = $t.Cell(13, 2)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "92÷4="

$cell = $t.Cell(13, 3)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "44÷9="

$cell = $t.Cell(13, 4)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "73÷2="

$cell = $t.Cell(13, 5)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "95÷9="

$cell = $t.Cell(17, 1)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "24÷5="

$cell = $t.Cell(17, 2)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "14÷6="

$cell = $t.Cell(17, 3)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "56÷4="

$cell = $t.Cell(17, 4)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "75÷4="

$cell = $t.Cell(17, 5)
$r2 = $cell.Range
$r2.SetRange($r2.Start, $r2.End - 2)
$r2.Text = "70÷5="

Write-Output "Applied date and table updates"